$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) — full cascade update
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Range("F2").Value = 1117
$ws1.Range("F3").Value = 824
$ws1.Range("F8").Value = 2046
$ws1.Range("F9").Value = 7641
$ws1.Range("F10").Value = 911
$ws1.Range("F11").Value = 430
$ws1.Range("F15").Value = 157
$ws1.Range("F16").Value = 7794
$ws1.Range("F17").Value = 310
$ws1.Range("F18").Value = 1353
$ws1.Range("G21").Value = "已售罄"
$ws1.Range("C22").Value = "江西·ShiningStaR数字互娱嘉年华配音演员史泽鲲专场见面会"
$ws1.Range("D22").Value = "前湖大道欣悦湖体育馆 欣悦湖体育馆"
$ws1.Range("E22").Value = "2024.05.03 09:30-05.03 17:30"
$ws1.Range("F22").Value = 228
$ws1.Range("G22").Value = "已售罄"
$ws1.Range("H22").Value = "https://show.bilibili.com/platform/detail.html?id=83497"
$ws1.Range("I22").Value = "//i1.hdslb.com/bfs/openplatform/202403/qm19B8RF1711620646864.jpeg"
$ws1.Range("C23").Value = "萍乡·AU8春季国漫展"
$ws1.Range("D23").Value = "金陵东路18号 萍乡市体育馆"
$ws1.Range("E23").Value = "2024.05.03 10:00-05.03 17:00"
$ws1.Range("F23").Value = 154
$ws1.Range("H23").Value = "https://show.bilibili.com/platform/detail.html?id=84259"
$ws1.Range("I23").Value = "//i0.hdslb.com/bfs/openplatform/202404/GxgrLKhd1712560102907.png"
$ws1.Range("B24").Value = "2024-05-03"
$ws1.Range("C24").Value = "赣州·漫库书店次元漫展"
$ws1.Range("D24").Value = "南门口地一大道下沉广场 漫库书店"
$ws1.Range("E24").Value = "2024.05.03 10:00-05.04 18:00"
$ws1.Range("F24").Value = 312
$ws1.Range("G24").Value = 45
$ws1.Range("H24").Value = "https://show.bilibili.com/platform/detail.html?id=83855"
$ws1.Range("I24").Value = "//i0.hdslb.com/bfs/openplatform/202404/juDVRy6Y1712481590113.jpeg"
$ws1.Range("C25").Value = "南昌·Youth动漫美食嘉年华"
$ws1.Range("D25").Value = "灌婴路西口朝阳江滩公园内 元亨体育运动中心"
$ws1.Range("E25").Value = "2024.05.04 09:00-05.05 20:00"
$ws1.Range("F25").Value = 146
$ws1.Range("G25").Value = 55
$ws1.Range("H25").Value = "https://show.bilibili.com/platform/detail.html?id=84036"
$ws1.Range("I25").Value = "//i2.hdslb.com/bfs/openplatform/202404/I5vd5js01712648032400.jpeg"
$ws1.Range("B26").Value = "2024-05-04"
$ws1.Range("C26").Value = "江西·ShiningStaR数字互娱嘉年华 配音演员陈张太康、张惠霖专场见面会"
$ws1.Range("D26").Value = "前湖大道欣悦湖体育馆 欣悦湖体育馆"
$ws1.Range("E26").Value = "2024.05.04 09:30-05.04 17:30"
$ws1.Range("F26").Value = 164
$ws1.Range("G26").Value = 228
$ws1.Range("H26").Value = "https://show.bilibili.com/platform/detail.html?id=83593"
$ws1.Range("I26").Value = "//i0.hdslb.com/bfs/openplatform/202404/LcnCzDxF1711935576170.jpeg"
$ws1.Range("B27").Value = "2024-05-12"
$ws1.Range("C27").Value = "宜春·BM次元盛典运动番only"
$ws1.Range("D27").Value = "鼓楼西路与官圳路交叉口东120米 地中海宴会酒店(润达店)"
$ws1.Range("E27").Value = "2024.05.12 10:00-05.12 17:00"
$ws1.Range("F27").Value = 20
$ws1.Range("G27").Value = 55
$ws1.Range("H27").Value = "https://show.bilibili.com/platform/detail.html?id=84636"
$ws1.Range("I27").Value = "//i2.hdslb.com/bfs/openplatform/202404/sNKPZWMh1713518729449.png"
$ws1.Range("C28").Value = "南昌·花绒万兽首届兽聚"
$ws1.Range("D28").Value = "南龙潘街666号二楼万达影城斜对面 融创茂啃趣馆"
$ws1.Range("E28").Value = "2024.05.18 09:30-05.19 16:30"
$ws1.Range("F28").Value = 110
$ws1.Range("H28").Value = "https://show.bilibili.com/platform/detail.html?id=83689"
$ws1.Range("I28").Value = "//i2.hdslb.com/bfs/openplatform/202403/h4iL6IvI1711790121140.jpeg"
$ws1.Range("B29").Value = "2024-05-18"
$ws1.Range("C29").Value = "赣州·原铁崩only"
$ws1.Range("D29").Value = "金岭东大道18号 万达广场(赣州经开店)"
$ws1.Range("E29").Value = "2024.05.18 10:00-05.19 17:00"
$ws1.Range("F29").Value = 23
$ws1.Range("G29").Value = 60
$ws1.Range("H29").Value = "https://show.bilibili.com/platform/detail.html?id=84721"
$ws1.Range("I29").Value = "//i0.hdslb.com/bfs/openplatform/202404/0n0MQiZh1713505673648.jpeg"
$ws1.Range("B30").Value = "2024-05-26"
$ws1.Range("C30").Value = "南昌·代号鸢盛花行only"
$ws1.Range("D30").Value = "民德路411号 东方豪景花园酒店(民德路店)"
$ws1.Range("E30").Value = "2024.05.26 09:30-05.26 17:30"
$ws1.Range("F30").Value = 416
$ws1.Range("G30").Value = 78
$ws1.Range("H30").Value = "https://show.bilibili.com/platform/detail.html?id=82529"
$ws1.Range("I30").Value = "//i2.hdslb.com/bfs/openplatform/202404/talOodLW1714030986517.png"
$ws1.Range("B31").Value = "2024-06-01"
$ws1.Range("C31").Value = "南昌·ACG CLUB动漫游戏嘉年华"
$ws1.Range("D31").Value = "火炬五路869号(科技城地铁站3号口步行340米) Ai羽球馆"
$ws1.Range("E31").Value = "2024.06.01 10:00-06.01 17:00"
$ws1.Range("F31").Value = 1133
$ws1.Range("G31").Value = 55
$ws1.Range("H31").Value = "https://show.bilibili.com/platform/detail.html?id=84497"
$ws1.Range("I31").Value = "//i1.hdslb.com/bfs/openplatform/202404/hZdMDMTZ1713768751631.jpeg"
$ws1.Range("B32").Value = "2024-06-09"
$ws1.Range("C32").Value = "信丰·端午节UPUP动漫展"
$ws1.Range("D32").Value = "迎宾大道富华双钻名汇西南侧约200米 诚瑞橙子体育馆"
$ws1.Range("E32").Value = "2024.06.09 10:00-06.09 17:00"
$ws1.Range("F32").Value = 55
$ws1.Range("G32").Value = 48
$ws1.Range("H32").Value = "https://show.bilibili.com/platform/detail.html?id=84078"
$ws1.Range("I32").Value = "//i0.hdslb.com/bfs/openplatform/202404/Qy0EOl551712651477492.jpeg"
$ws1.Range("F33").Value = 62
$ws1.Range("F34").Value = 80
$ws1.Range("F35").Value = 40

# Sheet "全部类型" (All types) — counter-only update
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("F2").Value = 1117
$ws4.Range("F3").Value = 824
$ws4.Range("F8").Value = 2046
$ws4.Range("F9").Value = 7641
$ws4.Range("F10").Value = 911
$ws4.Range("F11").Value = 430
$ws4.Range("F15").Value = 157
$ws4.Range("F16").Value = 7794
$ws4.Range("F17").Value = 310
$ws4.Range("F18").Value = 1353
$ws4.Range("G21").Value = "已售罄"
$ws4.Range("F27").Value = 110
$ws4.Range("F29").Value = 416
$ws4.Range("F30").Value = 1133
$ws4.Range("F33").Value = 62
$ws4.Range("F34").Value = 80
$ws4.Range("F35").Value = 40

